# cambios de mayo de mayo
# Update the reporting-period row (row 8) from the Q4-2021 period to the
# Q1-2022 period, and move the active selection/scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ejercicio: 2021 -> 2022
$ws.Range("A8").Value = 2022

# Fecha de inicio del periodo que se informa: 2021-10-01 -> 2022-01-01
$ws.Range("B8").Value = 44562

# Fecha de término del periodo que se informa: 2021-12-31 -> 2022-03-31
$ws.Range("C8").Value = 44651

# Fecha de elaboración: 2022-01-10 -> 2022-04-08
$ws.Range("E8").Value = 44659

# Fecha de validación: 2022-01-10 -> 2022-04-08
$ws.Range("H8").Value = 44659

# Fecha de actualización: 2022-01-10 -> 2022-04-08
$ws.Range("I8").Value = 44659

# Move the selection (was J8) to D12
$ws.Range("D12").Select()
